$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 95
$ws.Range("B2").Value = 2
$ws.Range("B3").Value = 9
$ws.Range("B4").Value = 86
$ws.Range("B5").Value = 41
$ws.Range("B6").Value = 66
$ws.Range("B7").Value = 63
$ws.Range("B8").Value = 24
$ws.Range("B9").Value = 69
$ws.Range("B10").Value = 40
